# Update Cavezzo.xlsx: append 4 new daily rows (230-233) after the
# existing last row (229), following the same A:D layout and carrying
# column A's date style ("s=2") down onto the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (229) down into the
# four new rows so the new A230:A233 cells inherit the date cell style.
$ws.Range("A229:D229").Copy()
$ws.Range("A230:D233").PasteSpecial(-4122)  # xlPasteFormats

$newRows = @(
    @(44304, 3, 12, 170.697012802276),
    @(44305, 1, 11, 156.4722617354196),
    @(44306, 0, 10, 142.2475106685633),
    @(44307, 0, 10, 142.2475106685633)
)

$r = 230
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

Write-Output "Added rows 230-233 to Sheet1"
